$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion note text (A1) ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 3.5 = 13340.75 pesos`n✅ 13340.75 pesos = 3.49 = 952.16 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $newText

# --- Sheet "tasas": update the rate cells N10, O10, N12, O12 ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 285.849
$wsTasas.Range("O10").Value = 3813.44
$wsTasas.Range("N12").Value = 3825
$wsTasas.Range("O12").Value = 273
